# Weekly Fruta/Hortaliza price update.
# A new week of price data (2023-04-25, serial 45041) was added at the top
# of the "Vega Monumental Concepción - Pera" data block, pushing the
# existing rows down by two. Concretely: insert two new rows at 635/636
# and populate them with the new week's "Packham's Triumph" Primera /
# Segunda records; everything below shifts down automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 635 (inherits the
# formatting -- including the date NumberFormat on column D -- from the
# row above, same as interactive Excel "Insert Copied/Sheet Rows").
$ws.Rows("635:636").Insert()

# New row 635: Packham's Triumph, Primera
$ws.Range("A635").Value = 11
$ws.Range("B635").Value = "Vega Monumental Concepción"
$ws.Range("C635").Value = "Bíobío"
$ws.Range("D635").Value = 45041
$ws.Range("E635").Value = 8
$ws.Range("F635").Value = "Fruta"
$ws.Range("G635").Value = 100104
$ws.Range("H635").Value = "Frutos de pepita"
$ws.Range("I635").Value = 100104005
$ws.Range("J635").Value = "Pera"
$ws.Range("K635").Value = "Packham's Triumph"
$ws.Range("L635").Value = "Primera"
$ws.Range("M635").Value = 100
$ws.Range("N635").Value = 10000
$ws.Range("O635").Value = 11000
$ws.Range("P635").Value = 10500
$ws.Range("Q635").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R635").Value = "Región de O'Higgins"
$ws.Range("S635").Value = 656
$ws.Range("T635").Value = 16

# New row 636: Packham's Triumph, Segunda
$ws.Range("A636").Value = 11
$ws.Range("B636").Value = "Vega Monumental Concepción"
$ws.Range("C636").Value = "Bíobío"
$ws.Range("D636").Value = 45041
$ws.Range("E636").Value = 8
$ws.Range("F636").Value = "Fruta"
$ws.Range("G636").Value = 100104
$ws.Range("H636").Value = "Frutos de pepita"
$ws.Range("I636").Value = 100104005
$ws.Range("J636").Value = "Pera"
$ws.Range("K636").Value = "Packham's Triumph"
$ws.Range("L636").Value = "Segunda"
$ws.Range("M636").Value = 50
$ws.Range("N636").Value = 9000
$ws.Range("O636").Value = 9000
$ws.Range("P636").Value = 9000
$ws.Range("Q636").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R636").Value = "Región de O'Higgins"
$ws.Range("S636").Value = 562
$ws.Range("T636").Value = 16
